$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(3, 8).Value = 6531.75
$ws.Cells.Item(3, 10).Value = 7416.6665
$ws.Cells.Item(3, 12).Value = 7416.6665
$ws.Cells.Item(3, 14).Value = -7644.6665
$ws.Cells.Item(12, 8).Value = 25000446
$ws.Cells.Item(12, 9).Value = 574.3333
$ws.Cells.Item(12, 11).Value = 574.3333
$ws.Cells.Item(12, 13).Value = -404.3333
$ws.Cells.Item(19, 8).Value = 1122.2307
$ws.Cells.Item(19, 9).Value = 971.1429000000001
$ws.Cells.Item(19, 10).Value = 1298.5
$ws.Cells.Item(19, 11).Value = 971.1429000000001
$ws.Cells.Item(19, 12).Value = 1298.5
$ws.Cells.Item(19, 13).Value = -796.1429000000001
$ws.Cells.Item(19, 14).Value = -1648.5
$ws.Cells.Item(33, 8).Value = 1068.9166
$ws.Cells.Item(33, 10).Value = 1499.5
$ws.Cells.Item(33, 12).Value = 1499.5
$ws.Cells.Item(33, 14).Value = -1957.5
$ws.Cells.Item(53, 8).Value = 30
$ws.Cells.Item(53, 9).Value = 33.6
$ws.Cells.Item(53, 10).Value = 27.428572
$ws.Cells.Item(53, 11).Value = 33.6
$ws.Cells.Item(53, 12).Value = 27.428572
$ws.Cells.Item(53, 13).Value = 603.4
$ws.Cells.Item(53, 14).Value = -1301.428572
$ws.Cells.Item(64, 8).Value = 4999.4546
$ws.Cells.Item(67, 8).Value = 4999.4546
$ws.Cells.Item(102, 8).Value = 6531.75
$ws.Cells.Item(102, 10).Value = 7416.6665
$ws.Cells.Item(102, 12).Value = 7416.6665
$ws.Cells.Item(102, 14).Value = -13906.6665
$ws.Cells.Item(108, 8).Value = 66657.664
$ws.Cells.Item(108, 10).Value = 66657.664
$ws.Cells.Item(108, 12).Value = 66657.664
$ws.Cells.Item(108, 14).Value = -74337.664
$ws.Cells.Item(109, 8).Value = 94608.336
$ws.Cells.Item(109, 10).Value = 94608.336
$ws.Cells.Item(109, 12).Value = 94608.336
$ws.Cells.Item(109, 14).Value = -97382.336
$ws.Cells.Item(110, 8).Value = 59534.4
$ws.Cells.Item(110, 10).Value = 59534.4
$ws.Cells.Item(110, 12).Value = 59534.4
$ws.Cells.Item(110, 14).Value = -67714.39999999999
$ws.Cells.Item(112, 8).Value = 1250.8636
$ws.Cells.Item(112, 9).Value = 2500
$ws.Cells.Item(112, 11).Value = 7500
$ws.Cells.Item(112, 13).Value = -6392
$ws.Cells.Item(117, 8).Value = 91321.73
$ws.Cells.Item(117, 10).Value = 91321.73
$ws.Cells.Item(117, 12).Value = 91321.73
$ws.Cells.Item(117, 14).Value = -100499.73
$ws.Cells.Item(123, 8).Value = 84670
$ws.Cells.Item(123, 10).Value = 84670
$ws.Cells.Item(123, 12).Value = 84670
$ws.Cells.Item(123, 14).Value = -94470
$ws.Cells.Item(135, 8).Value = 1041.1052
$ws.Cells.Item(135, 9).Value = 1093.3889
$ws.Cells.Item(135, 11).Value = 9840.500099999999
$ws.Cells.Item(135, 13).Value = -7305.500099999999
$ws.Cells.Item(138, 8).Value = 2387.0977
$ws.Cells.Item(138, 10).Value = 3234.739
$ws.Cells.Item(138, 12).Value = 9704.217000000001
$ws.Cells.Item(138, 14).Value = -19984.217

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(4, 8).Value = 500
$ws.Cells.Item(4, 10).Value = 0
$ws.Cells.Item(4, 12).Value = 0
$ws.Cells.Item(4, 14).ClearContents()
$ws.Cells.Item(45, 8).Value = 3535
$ws.Cells.Item(45, 9).Value = 3496.6667
$ws.Cells.Item(45, 11).Value = 3496.6667
$ws.Cells.Item(45, 13).Value = -3119.6667
$ws.Cells.Item(74, 8).Value = 49587.43
$ws.Cells.Item(74, 9).Value = 84469.586
$ws.Cells.Item(74, 11).Value = 84469.586
$ws.Cells.Item(74, 13).Value = -83595.586
$ws.Cells.Item(77, 8).Value = 49587.43
$ws.Cells.Item(77, 9).Value = 84469.586
$ws.Cells.Item(77, 11).Value = 422347.93
$ws.Cells.Item(77, 13).Value = -417979.93
$ws.Cells.Item(110, 8).Value = 1299.4375
$ws.Cells.Item(110, 9).Value = 1237.9333
$ws.Cells.Item(110, 11).Value = 1237.9333
$ws.Cells.Item(110, 13).Value = 807.0667000000001
$ws.Cells.Item(122, 8).Value = 9526.25
$ws.Cells.Item(122, 9).Value = 9526.25
$ws.Cells.Item(122, 11).Value = 28578.75
$ws.Cells.Item(122, 13).Value = -26128.75
$ws.Cells.Item(132, 8).Value = 1845.75
$ws.Cells.Item(132, 9).Value = 1628.5
$ws.Cells.Item(132, 10).Value = 2714.75
$ws.Cells.Item(132, 11).Value = 4885.5
$ws.Cells.Item(132, 12).Value = 8144.25
$ws.Cells.Item(132, 13).Value = -2355.5
$ws.Cells.Item(132, 14).Value = -13204.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(22, 8).Value = 104662.4
$ws.Cells.Item(22, 9).Value = 167680.17
$ws.Cells.Item(22, 10).Value = 10135.75
$ws.Cells.Item(22, 11).Value = 167680.17
$ws.Cells.Item(22, 12).Value = 10135.75
$ws.Cells.Item(22, 13).Value = -167507.17
$ws.Cells.Item(22, 14).Value = -10481.75
$ws.Cells.Item(99, 8).Value = 3908456
$ws.Cells.Item(99, 9).Value = 2274.8333
$ws.Cells.Item(99, 10).Value = 15627000
$ws.Cells.Item(99, 11).Value = 2274.8333
$ws.Cells.Item(99, 12).Value = 15627000
$ws.Cells.Item(99, 13).Value = -776.8332999999998
$ws.Cells.Item(99, 14).Value = -15629996
$ws.Cells.Item(134, 8).Value = 2040.8096
$ws.Cells.Item(134, 9).Value = 1337.4667
$ws.Cells.Item(134, 11).Value = 4012.4001
$ws.Cells.Item(134, 13).Value = -1477.4001

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 14610.781
$ws.Cells.Item(31, 9).Value = 3275.2222
$ws.Cells.Item(31, 10).Value = 19046.436
$ws.Cells.Item(31, 11).Value = 3275.2222
$ws.Cells.Item(31, 12).Value = 19046.436
$ws.Cells.Item(31, 13).Value = -2980.2222
$ws.Cells.Item(31, 14).Value = -19636.436
$ws.Cells.Item(34, 8).Value = 14610.781
$ws.Cells.Item(34, 9).Value = 3275.2222
$ws.Cells.Item(34, 10).Value = 19046.436
$ws.Cells.Item(34, 11).Value = 3275.2222
$ws.Cells.Item(34, 12).Value = 19046.436
$ws.Cells.Item(34, 13).Value = -3073.2222
$ws.Cells.Item(34, 14).Value = -19450.436
$ws.Cells.Item(86, 8).Value = 4682.5
$ws.Cells.Item(86, 9).Value = 4723.75
$ws.Cells.Item(86, 10).Value = 4600
$ws.Cells.Item(86, 11).Value = 4723.75
$ws.Cells.Item(86, 12).Value = 4600
$ws.Cells.Item(86, 13).Value = -3600.75
$ws.Cells.Item(86, 14).Value = -6846
$ws.Cells.Item(89, 8).Value = 4682.5
$ws.Cells.Item(89, 9).Value = 4723.75
$ws.Cells.Item(89, 10).Value = 4600
$ws.Cells.Item(89, 11).Value = 23618.75
$ws.Cells.Item(89, 12).Value = 23000
$ws.Cells.Item(89, 13).Value = -18002.75
$ws.Cells.Item(89, 14).Value = -34232
$ws.Cells.Item(122, 8).Value = 4747.5713
$ws.Cells.Item(122, 9).Value = 4311
$ws.Cells.Item(122, 11).Value = 12933
$ws.Cells.Item(122, 13).Value = -10483

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(33, 8).Value = 164.88889
$ws.Cells.Item(33, 9).Value = 92.14286
$ws.Cells.Item(33, 11).Value = 552.85716
$ws.Cells.Item(33, 13).Value = -269.85716
$ws.Cells.Item(107, 8).Value = 377.7143
$ws.Cells.Item(107, 9).Value = 365.33334
$ws.Cells.Item(107, 10).Value = 387
$ws.Cells.Item(107, 11).Value = 1096.00002
$ws.Cells.Item(107, 12).Value = 1161
$ws.Cells.Item(107, 13).Value = 823.9999800000001
$ws.Cells.Item(107, 14).Value = -5001

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(97, 8).Value = 526813.7
$ws.Cells.Item(97, 9).Value = 769573.6
$ws.Cells.Item(97, 10).Value = 833.8333
$ws.Cells.Item(97, 11).Value = 769573.6
$ws.Cells.Item(97, 12).Value = 833.8333
$ws.Cells.Item(97, 13).Value = -769077.6
$ws.Cells.Item(97, 14).Value = -1825.8333
$ws.Cells.Item(122, 8).Value = 18236.625
$ws.Cells.Item(122, 9).Value = 15883.8
$ws.Cells.Item(122, 11).Value = 47651.39999999999
$ws.Cells.Item(122, 13).Value = -45201.39999999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(40, 8).Value = 3971074.2
$ws.Cells.Item(40, 9).Value = 2406.2
$ws.Cells.Item(40, 10).Value = 13892744
$ws.Cells.Item(40, 11).Value = 2406.2
$ws.Cells.Item(40, 12).Value = 13892744
$ws.Cells.Item(40, 13).Value = -2270.2
$ws.Cells.Item(40, 14).Value = -13893016
$ws.Cells.Item(55, 8).Value = 1773.7858
$ws.Cells.Item(55, 9).Value = 1067.36
$ws.Cells.Item(55, 10).Value = 2812.647
$ws.Cells.Item(55, 11).Value = 1067.36
$ws.Cells.Item(55, 12).Value = 2812.647
$ws.Cells.Item(55, 13).Value = -894.3599999999999
$ws.Cells.Item(55, 14).Value = -3158.647
$ws.Cells.Item(93, 8).Value = 1043.2222
$ws.Cells.Item(93, 9).Value = 1168.75
$ws.Cells.Item(93, 10).Value = 942.8
$ws.Cells.Item(93, 11).Value = 1168.75
$ws.Cells.Item(93, 12).Value = 942.8
$ws.Cells.Item(93, 13).Value = 79.25
$ws.Cells.Item(93, 14).Value = -3438.8
$ws.Cells.Item(106, 8).Value = 34638.57
$ws.Cells.Item(106, 10).Value = 34638.57
$ws.Cells.Item(106, 12).Value = 34638.57
$ws.Cells.Item(106, 14).Value = -37162.57
$ws.Cells.Item(122, 8).Value = 20042462
$ws.Cells.Item(122, 9).Value = 59087.57
$ws.Cells.Item(122, 11).Value = 177262.71
$ws.Cells.Item(122, 13).Value = -174812.71
$ws.Cells.Item(125, 8).Value = 134999
$ws.Cells.Item(125, 10).Value = 134999
$ws.Cells.Item(125, 12).Value = 134999
$ws.Cells.Item(125, 14).Value = -144839
$ws.Cells.Item(136, 8).Value = 10314.143
$ws.Cells.Item(136, 9).Value = 5365.6665
$ws.Cells.Item(136, 11).Value = 16096.9995
$ws.Cells.Item(136, 13).Value = -13546.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(100, 8).Value = 4203107.5
$ws.Cells.Item(100, 9).Value = 7938400
$ws.Cells.Item(100, 11).Value = 15876800
$ws.Cells.Item(100, 13).Value = -15876259
$ws.Cells.Item(107, 8).Value = 7121.75
$ws.Cells.Item(107, 9).Value = 4731.357
$ws.Cells.Item(107, 10).Value = 8642.909
$ws.Cells.Item(107, 11).Value = 14194.071
$ws.Cells.Item(107, 12).Value = 25928.727
$ws.Cells.Item(107, 13).Value = -12274.071
$ws.Cells.Item(107, 14).Value = -29768.727
$ws.Cells.Item(116, 8).Value = 99999
$ws.Cells.Item(116, 10).Value = 99999
$ws.Cells.Item(116, 12).Value = 99999
$ws.Cells.Item(116, 14).Value = -109177
$ws.Cells.Item(122, 8).Value = 2650.9524
$ws.Cells.Item(122, 9).Value = 2436.0667
$ws.Cells.Item(122, 10).Value = 3188.1667
$ws.Cells.Item(122, 11).Value = 7308.2001
$ws.Cells.Item(122, 12).Value = 9564.500100000001
$ws.Cells.Item(122, 13).Value = -4858.2001
$ws.Cells.Item(122, 14).Value = -14464.5001
$ws.Cells.Item(125, 8).Value = 50997
$ws.Cells.Item(125, 10).Value = 50997
$ws.Cells.Item(125, 12).Value = 50997
$ws.Cells.Item(125, 14).Value = -60837
$ws.Cells.Item(132, 8).Value = 4832857
$ws.Cells.Item(132, 9).Value = 1690.4
$ws.Cells.Item(132, 11).Value = 5071.200000000001
$ws.Cells.Item(132, 13).Value = -2541.200000000001
